$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the two runs describing the phone number formats into a
#    single run/sentence: "Phone numbers will be accepted in three
#    formats and will be converted to +27:"
# ---------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "Phone numbers will be accepted in three formats and will be converted to +27:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Phone numbers will be accepted in three formats and will be converted to +27:",
    2)

# ---------------------------------------------------------------------
# 2. Mark the screenshot just above "Creating the VCF file" as
#    NoProof (adds <w:rPr><w:noProof/></w:rPr> to its run).
#    It is the second inline picture in the document (223.5pt x
#    181.5pt == 2838846 x 2305372 EMU).
# ---------------------------------------------------------------------
$shp = $d.InlineShapes.Item(2)
$shp.Range.NoProofing = 1

# ---------------------------------------------------------------------
# 3. Append " (if exe was created)" to the "Creating the VCF file"
#    heading as its own separate run.
# ---------------------------------------------------------------------
$rng = $d.Content
$found2 = $rng.Find.Execute("Creating the VCF file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" (if exe was created)")
# Toggle a character property on the freshly inserted text so the
# engine keeps it as its own run instead of re-coalescing it with the
# preceding "Creating the VCF file" run.
$rng.Font.Bold = 1
$rng.Font.Bold = 0
